$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value2 = $value
    $range.Style = "Normal"
}

# Row 2
$ws.Range("B2").Value2 = "Bitcoin"
$ws.Range("C2").Value2 = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
Set-TextValue $ws.Range("D2") "28.081.30"
$ws.Range("E2").Value2 = "  -4.73%  "

# Row 3
$ws.Range("B3").Value2 = "Ethereum"
$ws.Range("C3").Value2 = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
Set-TextValue $ws.Range("D3") "1.830.15"
$ws.Range("E3").Value2 = "  -3.51%  "

# Row 4
$ws.Range("B4").Value2 = "TetherUSD"
$ws.Range("C4").Value2 = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
Set-TextValue $ws.Range("D4") "0.9999"
$ws.Range("E4").Value2 = "  -0.47%  "

# Row 5
$ws.Range("B5").Value2 = "BNB"
$ws.Range("C5").Value2 = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
Set-TextValue $ws.Range("D5") "328.98"
$ws.Range("E5").Value2 = "  -2.90%  "

# Row 6
$ws.Range("B6").Value2 = "USDC"
$ws.Range("C6").Value2 = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue $ws.Range("D6") "0.9994"
$ws.Range("E6").Value2 = "  -0.35%  "

# Row 7
$ws.Range("B7").Value2 = "XRP"
$ws.Range("C7").Value2 = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextValue $ws.Range("D7") "0.4650"
$ws.Range("E7").Value2 = "  -2.20%  "

# Row 8
$ws.Range("B8").Value2 = "Cardano"
$ws.Range("C8").Value2 = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextValue $ws.Range("D8") "0.3867"
$ws.Range("E8").Value2 = "  -3.42%  "

# Row 9
$ws.Range("B9").Value2 = "Dogecoin"
$ws.Range("C9").Value2 = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue $ws.Range("D9") "0.07875"
$ws.Range("E9").Value2 = "  -1.91%  "

# Row 10
$ws.Range("B10").Value2 = "Polygon"
$ws.Range("C10").Value2 = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D10") "0.9604"
$ws.Range("E10").Value2 = "  -3.08%  "

# Row 11
$ws.Range("B11").Value2 = "Solana"
$ws.Range("C11").Value2 = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue $ws.Range("D11") "22.04"
$ws.Range("E11").Value2 = "  -5.00%  "

# Row 12
$ws.Range("B12").Value2 = "WrappedEther"
$ws.Range("C12").Value2 = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D12") "1.867.34"
$ws.Range("E12").Value2 = "  -3.13%  "

# Row 13
$ws.Range("B13").Value2 = "Polkadot"
$ws.Range("C13").Value2 = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D13") "5.660"
$ws.Range("E13").Value2 = "  -4.64%  "

# Row 14
$ws.Range("B14").Value2 = "Chainlink"
$ws.Range("C14").Value2 = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D14") "6.895"
$ws.Range("E14").Value2 = "  -2.68%  "

# Row 15
$ws.Range("B15").Value2 = "TRON"
$ws.Range("C15").Value2 = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D15") "0.06849"
$ws.Range("E15").Value2 = "  +0.56%  "

# Row 16
$ws.Range("B16").Value2 = "BinanceUSD"
$ws.Range("C16").Value2 = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue $ws.Range("D16") "0.9996"
$ws.Range("E16").Value2 = "  -0.54%  "

# Row 17
$ws.Range("B17").Value2 = "Litecoin"
$ws.Range("C17").Value2 = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D17") "86.66"
$ws.Range("E17").Value2 = "  -2.66%  "

# Row 18
$ws.Range("B18").Value2 = "ShibaInu"
$ws.Range("C18").Value2 = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D18") "0.000009999"
$ws.Range("E18").Value2 = "  -1.91%  "

# Row 19
$ws.Range("B19").Value2 = "Avalanche"
$ws.Range("C19").Value2 = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws.Range("D19") "16.66"
$ws.Range("E19").Value2 = "  -3.78%  "

# Row 20
$ws.Range("B20").Value2 = "Dai"
$ws.Range("C20").Value2 = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D20") "0.9993"
$ws.Range("E20").Value2 = "  -0.42%  "

# Row 21
$ws.Range("B21").Value2 = "WrappedBTC"
$ws.Range("C21").Value2 = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue $ws.Range("D21") "28.081.44"
$ws.Range("E21").Value2 = "  -4.73%  "

# Row 22
$ws.Range("B22").Value2 = "Uniswap"
$ws.Range("C22").Value2 = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D22") "5.335"
$ws.Range("E22").Value2 = "  -3.16%  "

# Row 23
$ws.Range("B23").Value2 = "Cosmos"
$ws.Range("C23").Value2 = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D23") "11.01"
$ws.Range("E23").Value2 = "  -5.43%  "

# Row 24
$ws.Range("B24").Value2 = "Toncoin"
$ws.Range("C24").Value2 = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D24") "2.095"
$ws.Range("E24").Value2 = "  -2.68%  "

# Row 25
$ws.Range("B25").Value2 = "WrappedliquidstakedEther2.0"
$ws.Range("C25").Value2 = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D25") "2.058.61"
$ws.Range("E25").Value2 = "  -4.51%  "

# Row 26
$ws.Range("B26").Value2 = "Monero"
$ws.Range("C26").Value2 = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D26") "152.32"
$ws.Range("E26").Value2 = "  -3.20%  "

# Row 27
$ws.Range("B27").Value2 = "EthereumClassic"
$ws.Range("C27").Value2 = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D27") "19.22"
$ws.Range("E27").Value2 = "  -2.22%  "

# Row 28
$ws.Range("B28").Value2 = "InternetComputer(DFINITY)"
$ws.Range("C28").Value2 = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D28") "5.784"
$ws.Range("E28").Value2 = "  -10.84%  "

# Row 29
$ws.Range("B29").Value2 = "LidoDAOToken"
$ws.Range("C29").Value2 = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D29") "1.976"
$ws.Range("E29").Value2 = "  -3.69%  "

# Row 30
$ws.Range("B30").Value2 = "BitcoinCash"
$ws.Range("C30").Value2 = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Range("D30") "117.07"
$ws.Range("E30").Value2 = "  -1.73%  "

# Row 31
$ws.Range("B31").Value2 = "ImmutableX"
$ws.Range("C31").Value2 = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D31") "0.9388"
$ws.Range("E31").Value2 = "  -5.60%  "

# Row 32
$ws.Range("B32").Value2 = "Stellar"
$ws.Range("C32").Value2 = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D32") "0.09265"
$ws.Range("E32").Value2 = "  -2.99%  "

# Row 33
$ws.Range("B33").Value2 = "Filecoin"
$ws.Range("C33").Value2 = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D33") "5.305"
$ws.Range("E33").Value2 = "  -3.24%  "

# Row 34
$ws.Range("B34").Value2 = "ARBITRUM"
$ws.Range("C34").Value2 = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D34") "1.320"
$ws.Range("E34").Value2 = "  -4.71%  "

# Row 35
$ws.Range("B35").Value2 = "HuobiToken"
$ws.Range("C35").Value2 = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D35") "3.345"
$ws.Range("E35").Value2 = "  -5.23%  "

# Row 36
$ws.Range("B36").Value2 = "Hedera"
$ws.Range("C36").Value2 = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D36") "0.05943"
$ws.Range("E36").Value2 = "  -7.18%  "

# Row 37
$ws.Range("B37").Value2 = "VeChain"
$ws.Range("C37").Value2 = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D37") "0.02145"
$ws.Range("E37").Value2 = "  -4.39%  "

# Row 38
$ws.Range("B38").Value2 = "TrustWalletToken"
$ws.Range("C38").Value2 = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D38") "1.148"
$ws.Range("E38").Value2 = "  -4.36%  "

# Row 39
$ws.Range("B39").Value2 = "Frax"
$ws.Range("C39").Value2 = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue $ws.Range("D39") "0.9990"
$ws.Range("E39").Value2 = "  -0.45%  "

# Row 40
$ws.Range("B40").Value2 = "FraxShare"
$ws.Range("C40").Value2 = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D40") "7.658"
$ws.Range("E40").Value2 = "  -1.11%  "

# Row 41
$ws.Range("B41").Value2 = "TheSandbox"
$ws.Range("C41").Value2 = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D41") "0.5599"
$ws.Range("E41").Value2 = "  -3.95%  "

# Row 42
$ws.Range("B42").Value2 = "Aptos"
$ws.Range("C42").Value2 = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D42") "9.931"
$ws.Range("E42").Value2 = "  -5.75%  "

# Row 43
$ws.Range("B43").Value2 = "Algorand"
$ws.Range("C43").Value2 = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D43") "0.1772"
$ws.Range("E43").Value2 = "  -2.62%  "

# Row 44
$ws.Range("B44").Value2 = "WEMIXToken"
$ws.Range("C44").Value2 = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D44") "1.224"
$ws.Range("E44").Value2 = "  -3.08%  "

# Row 45
$ws.Range("B45").Value2 = "RenderToken"
$ws.Range("C45").Value2 = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D45") "2.233"
$ws.Range("E45").Value2 = "  -7.53%  "

# Row 46
$ws.Range("B46").Value2 = "EnergySwap"
$ws.Range("C46").Value2 = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D46") "11.60"
$ws.Range("E46").Value2 = "  -4.80%  "

# Row 47
$ws.Range("B47").Value2 = "Decentraland"
$ws.Range("C47").Value2 = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Range("D47") "0.5278"
$ws.Range("E47").Value2 = "  -4.00%  "

# Row 48
$ws.Range("B48").Value2 = "Cronos"
$ws.Range("C48").Value2 = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D48") "0.07021"
$ws.Range("E48").Value2 = "  -4.32%  "

# Row 49
$ws.Range("B49").Value2 = "NEARProtocol"
$ws.Range("C49").Value2 = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D49") "1.832"
$ws.Range("E49").Value2 = "  -6.17%  "

# Row 50
$ws.Range("B50").Value2 = "Quant"
$ws.Range("C50").Value2 = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D50") "111.56"
$ws.Range("E50").Value2 = "  -4.25%  "

# Row 51
$ws.Range("B51").Value2 = "PaxDollar"
$ws.Range("C51").Value2 = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws.Range("D51") "0.9990"
$ws.Range("E51").Value2 = "  -0.54%  "
